# Request Form - F1.xlsx : "sops Update 3"
#
# - Rename the visible worksheet from "F-SW-FR-01" to "S-SW-SC-01"
#   (this is a Software Service Catalog form, not a New Feature Request
#   form, so the sheet/code was relabelled accordingly).
# - Keep the printed area pinned to A1:F22 on the renamed sheet (the
#   print-area defined name is sheet-name-qualified, so it must be
#   re-pointed at the new sheet name).
# - Move the active selection from C10:F10 to A4:F4 (cursor left where
#   the editor was last working).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (fall back to a name lookup in case ActiveSheet isn't
# already the form sheet).
if ($ws.Name -ne "F-SW-FR-01") {
    $ws = $wb.Worksheets.Item("F-SW-FR-01")
}
$ws.Name = "S-SW-SC-01"

# Re-establish the print area under the new sheet name.
$ws.PageSetup.PrintArea = '$A$1:$F$22'

# Activate the sheet and move the selection to A4:F4.
$ws.Activate()
$ws.Range("A4:F4").Select()
